# Scheduled-runner refresh of the market-price / profit columns (H:N) on each
# job sheet. Values are re-pulled from the price source; only the H-N cell
# contents move (no formulas, no structural changes).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 321
$ws.Range("I2").Value = 533.3333
$ws.Range("J2").Value = 161.75
$ws.Range("K2").Value = 533.3333
$ws.Range("L2").Value = 161.75
$ws.Range("M2").Value = -420.3333
$ws.Range("N2").Value = -387.75
# row 17
$ws.Range("H17").Value = 1234.0667
$ws.Range("J17").Value = 1010.6818
$ws.Range("L17").Value = 3032.0454
$ws.Range("N17").Value = -3368.0454
# row 135
$ws.Range("H135").Value = 524.82355
$ws.Range("I135").Value = 506.76923
$ws.Range("K135").Value = 4560.92307
$ws.Range("M135").Value = -2025.92307
# row 137
$ws.Range("H137").Value = 1525.0714
$ws.Range("I137").Value = 1362.5834
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 4087.7502
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -1537.7502
$ws.Range("N137").Value = -12600
# row 138
$ws.Range("H138").Value = 2884.4888
$ws.Range("I138").Value = 2746.52
$ws.Range("J138").Value = 3056.95
$ws.Range("K138").Value = 8239.559999999999
$ws.Range("L138").Value = 9170.849999999999
$ws.Range("M138").Value = -3099.559999999999
$ws.Range("N138").Value = -19450.85

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3603.8723
$ws.Range("I32").Value = 2590.7805
$ws.Range("J32").Value = 10526.667
$ws.Range("K32").Value = 2590.7805
$ws.Range("L32").Value = 10526.667
$ws.Range("M32").Value = -2303.7805
$ws.Range("N32").Value = -11100.667
# row 45
$ws.Range("H45").Value = 1596.5333
$ws.Range("I45").Value = 1279.5714
$ws.Range("J45").Value = 1873.875
$ws.Range("K45").Value = 1279.5714
$ws.Range("L45").Value = 1873.875
$ws.Range("M45").Value = -902.5714
$ws.Range("N45").Value = -2627.875
# row 61
$ws.Range("H61").Value = 3752.4644
$ws.Range("I61").Value = 2721.4
$ws.Range("J61").Value = 12344.667
$ws.Range("K61").Value = 2721.4
$ws.Range("L61").Value = 12344.667
$ws.Range("M61").Value = -2509.4
$ws.Range("N61").Value = -12768.667
# row 74
$ws.Range("H74").Value = 1302.25
$ws.Range("I74").Value = 859.08
$ws.Range("K74").Value = 859.08
$ws.Range("M74").Value = 14.91999999999996
# row 77
$ws.Range("H77").Value = 1302.25
$ws.Range("I77").Value = 859.08
$ws.Range("K77").Value = 4295.400000000001
$ws.Range("M77").Value = 72.59999999999945
# row 136
$ws.Range("H136").Value = 3752.4644
$ws.Range("I136").Value = 2721.4
$ws.Range("J136").Value = 12344.667
$ws.Range("K136").Value = 8164.200000000001
$ws.Range("L136").Value = 37034.001
$ws.Range("M136").Value = -5614.200000000001
$ws.Range("N136").Value = -42134.001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 81879.75999999999
$ws.Range("I86").Value = 1720.7894
$ws.Range("K86").Value = 1720.7894
$ws.Range("M86").Value = -597.7893999999999
# row 89
$ws.Range("H89").Value = 81879.75999999999
$ws.Range("I89").Value = 1720.7894
$ws.Range("K89").Value = 8603.947
$ws.Range("M89").Value = -2987.947
# row 105
$ws.Range("H105").Value = 2555.5557
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 1433.1666
$ws.Range("I22").Value = 299.5
$ws.Range("K22").Value = 299.5
$ws.Range("M22").Value = 50.5
# row 31
$ws.Range("H31").Value = 2452.7273
$ws.Range("I31").Value = 1608.0322
$ws.Range("K31").Value = 1608.0322
$ws.Range("M31").Value = -1313.0322
# row 34
$ws.Range("H34").Value = 2452.7273
$ws.Range("I34").Value = 1608.0322
$ws.Range("K34").Value = 1608.0322
$ws.Range("M34").Value = -1406.0322
# row 58
$ws.Range("H58").Value = 1611855.5
$ws.Range("I58").Value = 3624267.2
$ws.Range("K58").Value = 3624267.2
$ws.Range("M58").Value = -3624064.2
# row 132
$ws.Range("H132").Value = 2600.2856
$ws.Range("I132").Value = 1556.1111
$ws.Range("K132").Value = 4668.3333
$ws.Range("M132").Value = -2138.3333
# row 134
$ws.Range("H134").Value = 1007.6667
$ws.Range("I134").Value = 1008.36365
$ws.Range("K134").Value = 3025.09095
$ws.Range("M134").Value = -490.0909499999998
# row 136
$ws.Range("H136").Value = 1611855.5
$ws.Range("I136").Value = 3624267.2
$ws.Range("K136").Value = 10872801.6
$ws.Range("M136").Value = -10870251.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 187.14285
$ws.Range("I2").Value = 320
$ws.Range("J2").Value = 87.5
$ws.Range("K2").Value = 1920
$ws.Range("L2").Value = 525
$ws.Range("M2").Value = -1807
$ws.Range("N2").Value = -751
# row 92
$ws.Range("H92").Value = 316.66666
$ws.Range("J92").Value = 340
$ws.Range("L92").Value = 1020
$ws.Range("N92").Value = -3516
# row 105
$ws.Range("H105").Value = 3080
$ws.Range("J105").Value = 3080
$ws.Range("L105").Value = 9240
$ws.Range("N105").Value = -14482
# row 107
$ws.Range("H107").Value = 809.5
$ws.Range("J107").Value = 809.5
$ws.Range("L107").Value = 2428.5
$ws.Range("N107").Value = -6268.5
# row 109
$ws.Range("H109").Value = 2196.7693
$ws.Range("I109").Value = 1009.75
$ws.Range("K109").Value = 3029.25
$ws.Range("M109").Value = -1989.25
# row 116
$ws.Range("H116").Value = 166668670
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
# row 121
$ws.Range("H121").Value = 641
$ws.Range("I121").Value = 382
$ws.Range("J121").Value = 900
$ws.Range("K121").Value = 1146
$ws.Range("L121").Value = 2700
$ws.Range("M121").Value = 164
$ws.Range("N121").Value = -5320
# row 129
$ws.Range("H129").Value = 49283.066
$ws.Range("J129").Value = 73456.10000000001
$ws.Range("L129").Value = 220368.3
$ws.Range("N129").Value = -230368.3
# row 131
$ws.Range("H131").Value = 9043.188
$ws.Range("I131").Value = 541.7778
$ws.Range("J131").Value = 10049.935
$ws.Range("K131").Value = 1625.3334
$ws.Range("L131").Value = 30149.805
$ws.Range("M131").Value = 3414.6666
$ws.Range("N131").Value = -40229.805
# row 137
$ws.Range("H137").Value = 3510
$ws.Range("I137").Value = 1303.3334
$ws.Range("J137").Value = 4111.8184
$ws.Range("K137").Value = 3910.0002
$ws.Range("L137").Value = 12335.4552
$ws.Range("M137").Value = 1189.9998
$ws.Range("N137").Value = -22535.4552

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# row 113
$ws.Range("H113").Value = 1300
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070
# row 122
$ws.Range("H122").Value = 1759.3077
$ws.Range("I122").Value = 1359
$ws.Range("J122").Value = 2226.3333
$ws.Range("K122").Value = 4077
$ws.Range("L122").Value = 6678.999899999999
$ws.Range("M122").Value = -1627
$ws.Range("N122").Value = -11578.9999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 3409.5
$ws.Range("I40").Value = 1152.8
$ws.Range("K40").Value = 1152.8
$ws.Range("M40").Value = -1016.8
# row 61
$ws.Range("H61").Value = 3112.75
$ws.Range("I61").Value = 2935.3
$ws.Range("K61").Value = 2935.3
$ws.Range("M61").Value = -2733.3
# row 113
$ws.Range("H113").Value = 3112.75
$ws.Range("I113").Value = 2935.3
$ws.Range("K113").Value = 2935.3
$ws.Range("M113").Value = -765.3000000000002
# row 132
$ws.Range("H132").Value = 2006.4445
$ws.Range("I132").Value = 1972.4546
$ws.Range("J132").Value = 2029.8125
$ws.Range("K132").Value = 5917.3638
$ws.Range("L132").Value = 6089.4375
$ws.Range("M132").Value = -3387.3638
$ws.Range("N132").Value = -11149.4375
# row 136
$ws.Range("H136").Value = 3709.8667
$ws.Range("I136").Value = 1664.1428
$ws.Range("J136").Value = 5499.875
$ws.Range("K136").Value = 4992.428400000001
$ws.Range("L136").Value = 16499.625
$ws.Range("M136").Value = -2442.428400000001
$ws.Range("N136").Value = -21599.625

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 1807.8055
$ws.Range("I132").Value = 1429.3334
$ws.Range("J132").Value = 2564.75
$ws.Range("K132").Value = 4288.0002
$ws.Range("L132").Value = 7694.25
$ws.Range("M132").Value = -1758.0002
$ws.Range("N132").Value = -12754.25
# row 136
$ws.Range("H136").Value = 15434285
$ws.Range("I136").Value = 21369416
$ws.Range("K136").Value = 64108248
$ws.Range("M136").Value = -64105698
